# Rename existing sheet "Sheet1" -> "mensual"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "mensual"

# Add a new worksheet "anual" right after "mensual"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "anual"

# Header row
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "tipo"
$ws2.Range("C1").Value = "cantidad"

# Annual summary data for 2023 and 2024
$data = @(
    @(2023, "Violencia doméstica (VD)", 872),
    @(2023, "Acecho (A)", 105),
    @(2023, "Agresión sexual (AS)", 25),
    @(2023, "Violencia en cita (VC)", 0),
    @(2023, "Discrimen de género (DG)", 92),
    @(2023, "Otras", 1797),
    @(2023, "Trata Humana", 1),
    @(2024, "Violencia doméstica (VD)", 1425),
    @(2024, "Acecho (A)", 85),
    @(2024, "Agresión sexual (AS)", 39),
    @(2024, "Violencia en cita (VC)", 2),
    @(2024, "Discrimen de género (DG)", 79),
    @(2024, "Otras", 2077),
    @(2024, "Trata Humana", 2)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Match the source's saved cursor positions on each sheet
[void]$ws2.Range("D18").Select()
[void]$ws1.Select()
[void]$ws1.Range("F34").Select()
